# Updated symbol list on Thu Dec 29 18:48:21 UTC 2022 with GitHub Actions
#
# Applies the refreshed "Price" (column D) and "Volume(1h)" (column E)
# figures scraped from coinranking.com. All of these cells are stored as
# plain text in the workbook (not numbers), so each write goes through a
# small helper that forces a text/string write while restoring the
# cell's original style afterwards (writing a numeric-looking string via
# .Value alone would silently coerce it to a Number cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$addr,
        [string]$text
    )

    $rng = $ws.Range($addr)
    $originalStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $originalStyle
}

# --- Column D: Price -------------------------------------------------
Set-TextValue "D2"  "245.79"
Set-TextValue "D3"  "24.21"
Set-TextValue "D4"  "5.283"
Set-TextValue "D6"  "6.500"
Set-TextValue "D7"  "3.146"
Set-TextValue "D8"  "0.8107"
Set-TextValue "D9"  "0.8639"
Set-TextValue "D10" "0.1379"
Set-TextValue "D11" "0.07000"
Set-TextValue "D12" "0.03153"
Set-TextValue "D13" "0.02913"
Set-TextValue "D14" "0.09386"
Set-TextValue "D15" "3.746"
Set-TextValue "D16" "0.001525"
Set-TextValue "D17" "0.04682"
Set-TextValue "D19" "0.006206"
Set-TextValue "D21" "0.004639"
Set-TextValue "D22" "0.00006102"
Set-TextValue "D23" "3.501"
Set-TextValue "D24" "2.148"
Set-TextValue "D40" "0.03712"
Set-TextValue "D41" "0.006282"
Set-TextValue "D42" "0.1053"
Set-TextValue "D43" "0.003201"
Set-TextValue "D44" "0.007750"
Set-TextValue "D45" "0.00005282"
Set-TextValue "D48" "0.002437"

# --- Column E: Volume(1h) label (rank+name+symbol+Best/Worst-in-24h) -
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
